$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; this shifts the existing rows 50..137
# down to 51..138, preserving all their values/styles, and leaves a
# blank row 50 (inheriting column D's date style) ready to populate.
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with this week's data point.
$ws.Range("A50").Value = 5
$ws.Range("B50").Value = "Macroferia Regional de Talca"
$ws.Range("C50").Value = "Maule"
$ws.Range("D50").Value = 44571
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = 100112031
$ws.Range("G50").Value = "Poroto verde"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("M50").Value = 30000
$ws.Range("N50").Value = "$/saco 25 kilos"
$ws.Range("O50").Value = "Región del Maule"
$ws.Range("P50").Value = 1200
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
